$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '26.715.92'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.600.71'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.15'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.0619'
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.64'
$ws.Range('E10').Value = '  +0.86%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').Value = '1.825.47'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '1.603.07'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('E14').Value = '  +0.49%  '
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.14'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '26.687.76'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '210.56'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.96'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.03'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.37'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0511'
$ws.Range('E30').Value = '  +0.65%  '
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D34').Value = '1.293.63'
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.48'
$ws.Range('E35').Value = '  +0.71%  '
$ws.Range('E36').Value = '  +1.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.606'
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('E38').Value = '  +17.14%  '
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.822'
$ws.Range('E40').Value = '  -1.83%  '
$ws.Range('E41').Value = '  -1.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.780'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.20'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = '1.738.64'
$ws.Range('E45').Value = '  +0.35%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '90.70'
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  -2.55%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.101'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0517'
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('E51').Value = '  -0.28%  '
